$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert the new "Data" worksheet between "About" and "CApULAbIFM"
# ------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsCApULAbIFM = $wb.Worksheets.Item("CApULAbIFM")
$wsData = $wb.Worksheets.Add($wsCApULAbIFM)
$wsData.Name = "Data"

# ------------------------------------------------------------------
# 2. "About" sheet - update source citation + add unit-conversion block
# ------------------------------------------------------------------
$wsAbout.Range("B3").Value = "Article: ""Indigenous Forests Are Some of the Amazon’s Last Carbon Sinks"""

$wsAbout.Range("B4").Value = "Veit, P., Gibbs, D., & Reytar, K. (2023). Indigenous Forests Are Some of the Amazon’s Last Carbon Sinks. World Resources Institute."
$wsAbout.Range("B4").WrapText = $true
$wsAbout.Rows.Item(4).RowHeight = 43.5

$wsAbout.Range("B5").Value = "https://www.wri.org/insights/amazon-carbon-sink-indigenous-forests#:~:text=Other%20Community%20Lands%20Are%20Also%20Strong%20Carbon%20Sinks&text=Our%20analysis%20of%20Afro%2Ddescendant,sinks%20from%202001%20to%202021.&text=Removals%20were%20about%20twice%20as,climate%20change%20through%20forest%20stewardship.&text=Collectively%20held%20forests%20in%20Mexico,the%20Philippines%20were%20carbon%20sinks."
$wsAbout.Hyperlinks.Add(
    $wsAbout.Range("B5"),
    "https://www.wri.org/insights/amazon-carbon-sink-indigenous-forests",
    ":~:text=Other%20Community%20Lands%20Are%20Also%20Strong%20Carbon%20Sinks&text=Our%20analysis%20of%20Afro%2Ddescendant,sinks%20from%202001%20to%202021.&text=Removals%20were%20about%20twice%20as,climate%20change%20through%20forest%20stewardship.&text=Collectively%20held%20forests%20in%20Mexico,the%20Philippines%20were%20carbon%20sinks.",
    "",
    "https://www.wri.org/insights/amazon-carbon-sink-indigenous-forests#:~:text=Other%20Community%20Lands%20Are%20Also%20Strong%20Carbon%20Sinks&text=Our%20analysis%20of%20Afro%2Ddescendant,sinks%20from%202001%20to%202021.&text=Removals%20were%20about%20twice%20as,climate%20change%20through%20forest%20stewardship.&text=Collectively%20held%20forests%20in%20Mexico,the%20Philippines%20were%20carbon%20sinks."
)

$wsAbout.Rows.Item(6).Delete()

$wsAbout.Range("A7").Value = "Unit Conversions"
$wsAbout.Range("A7").Font.Bold = $true
$wsAbout.Range("A7").WrapText = $true
$wsAbout.Rows.Item(7).RowHeight = 29

$wsAbout.Range("B7").Value = 2.471053815
$wsAbout.Range("C7").Value = "acre/ha"

$wsAbout.Range("B8").Value = 1000000
$wsAbout.Range("C8").Value = "g/ton"

$wsAbout.Columns.Item(1).ColumnWidth = 11.2734375
$wsAbout.Columns.Item(2).ColumnWidth = 41.21875

$wsAbout.Range("B8").Select()

# ------------------------------------------------------------------
# 3. "Data" sheet - new calculations
# ------------------------------------------------------------------
$wsData.Range("C2").Value = "2001 - 2021"

$wsData.Range("A3").Value = "Annual Tonnes of CO2 per hectare"
$wsData.Range("A3").WrapText = $true
$wsData.Rows.Item(3).RowHeight = 29
$wsData.Range("B3").Value = 1.65

$wsData.Range("A4").Value = "Annual grams of CO2 per acre"
$wsData.Range("A4").WrapText = $true
$wsData.Rows.Item(4).RowHeight = 29
$wsData.Range("B4").Formula = "=B3*About!B8/(About!B7)"
$wsData.Range("B4").NumberFormat = "0"

$wsData.Range("A6").Value = "Note: Figure given is for the mean of Indigenous and Afro-descent forests."

$wsData.Columns.Item(1).ColumnWidth = 19.2734375

$wsData.Range("C3").Select()

# ------------------------------------------------------------------
# 4. "CApULAbIFM" sheet - point the CO2 Abated figure at the Data sheet
# ------------------------------------------------------------------
$wsCApULAbIFM.Range("B2").Formula = "=Data!B4"

$wsCApULAbIFM.Columns.Item(1).ColumnWidth = 15.18359375
